# Insurance sheet ("保險", sheet 8): add company/insurance columns and
# reorder/extend data per commit "#5: insurance, claim, debt, investment done"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(8)

# --- Extend header row (row 1) formatting to the new columns E1:K1 ---
# Copy format (bold, centered, bordered) from the existing header cell D1
$ws.Range("D1").Copy($ws.Range("E1:K1"))

# --- Extend data rows (2-5) formatting to the new columns E:K ---
# Copy format (plain/default) from an existing data cell D2
$ws.Range("D2").Copy($ws.Range("E2:K5"))

# --- Header row values ---
$ws.Range("B1").Value2 = "company"
$ws.Range("C1").Value2 = "name"
$ws.Range("D1").Value2 = "owner"
$ws.Range("E1").Value2 = "property_category"
$ws.Range("F1").Value2 = "category"
$ws.Range("G1").Value2 = "date"
$ws.Range("H1").Value2 = "legislator_name"
$ws.Range("I1").Value2 = "legislator_id"
$ws.Range("J1").Value2 = "source_file"
$ws.Range("K1").Value2 = "index"

# --- Row 2 (index 134) ---
$ws.Range("B2").Value2 = "大都會國際人壽"
$ws.Range("C2").Value2 = "大都會終身壽險"
$ws.Range("D2").Value2 = "劉娟娟"
$ws.Range("E2").Value2 = "insurance"
$ws.Range("F2").Value2 = "normal"
$ws.Range("G2").NumberFormat = "@"
$ws.Range("G2").Value2 = "2013-11-22"
$ws.Range("H2").Value2 = "吳育昇"
$ws.Range("I2").Value2 = 1322
$ws.Range("J2").Value2 = "tmp88481"
$ws.Range("K2").Value2 = 134

# --- Row 3 (index 135) ---
$ws.Range("B3").Value2 = "富邦人壽"
$ws.Range("C3").Value2 = "富邦人壽子女儲蓄保險"
$ws.Range("D3").Value2 = "劉娟娟"
$ws.Range("E3").Value2 = "insurance"
$ws.Range("F3").Value2 = "normal"
$ws.Range("G3").NumberFormat = "@"
$ws.Range("G3").Value2 = "2013-11-22"
$ws.Range("H3").Value2 = "吳育昇"
$ws.Range("I3").Value2 = 1322
$ws.Range("J3").Value2 = "tmp88481"
$ws.Range("K3").Value2 = 135

# --- Row 4 (index 136) ---
$ws.Range("B4").Value2 = "南山人壽"
$ws.Range("C4").Value2 = "富邦人壽子女教育保險"
$ws.Range("D4").Value2 = "劉娟娟"
$ws.Range("E4").Value2 = "insurance"
$ws.Range("F4").Value2 = "normal"
$ws.Range("G4").NumberFormat = "@"
$ws.Range("G4").Value2 = "2013-11-22"
$ws.Range("H4").Value2 = "吳育昇"
$ws.Range("I4").Value2 = 1322
$ws.Range("J4").Value2 = "tmp88481"
$ws.Range("K4").Value2 = 136

# --- Row 5 (index 137) ---
$ws.Range("B5").Value2 = "南山人壽"
$ws.Range("C5").Value2 = "南山康寧終身壽險"
$ws.Range("D5").Value2 = "吳育昇"
$ws.Range("E5").Value2 = "insurance"
$ws.Range("F5").Value2 = "normal"
$ws.Range("G5").NumberFormat = "@"
$ws.Range("G5").Value2 = "2013-11-22"
$ws.Range("H5").Value2 = "吳育昇"
$ws.Range("I5").Value2 = 1322
$ws.Range("J5").Value2 = "tmp88481"
$ws.Range("K5").Value2 = 137

# Reset G2:G5 formatting back to the plain/default look (clears the "@" text
# format we used above to stop the date string from being auto-parsed as a
# serial date number, without touching the values already stored).
$ws.Range("D2").Copy()
$ws.Range("G2:G5").PasteSpecial(-4122)
